# Effort sheet: add a new row (24) for the newest log entry.
#   A24 = 2013-07-02 (serial 41457), formatted like the rows above it
#   B24 = 2.5 (Effort [h])
#   D24 = "Implementation tc14, variants tried" (new Task entry, appended
#         to the shared string table)
# Also moves the active-cell selection to C24, matching the saved view.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Effort R 1.0")

# Copy the date format (style index) from the row above so the new date
# cell keeps the existing "ddd dd/mm/yyyy" number format instead of
# creating a duplicate style entry.
$ws.Range("A23").Copy()
$ws.Range("A24").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Range("A24").Value = 41457
$ws.Range("B24").Value = 2.5
$ws.Range("D24").Value = "Implementation tc14, variants tried"

$ws.Range("C24").Select()
